$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row tweak ---
# B1: "Target OU (created if needed)" -> "Target OU (created if needed) (NEW EACH TIME)"
$ws.Range("B1").Value = "Target OU (created if needed) (NEW EACH TIME)"

# --- BES row (row 3): fill in the cart / OU name that had been left blank ---
# D3 gets a value and the yellow "needs attention" highlight is cleared.
$ws.Range("D3").Value = "BES Cart 4 2021 "
$ws.Range("D3").ClearFormats()
# C3 hyperlink text flips from the disabled "NOT..." address to the live address
# (the underlying mailto: target is left untouched).
$ws.Range("C3").Value = "ecarr@risd.k12.nm.us"

# --- NLE row (row 7): swap to the disabled contact and clear the cart name ---
$ws.Range("C7").Value = "NOTjsoltero@risd.k12.nm.us"
$ws.Range("D7").ClearContents()

# --- SMS row (row 18): the cart has been assigned already, clear the cart name ---
$ws.Range("D18").Value = ""

# --- Selection moves to E3 ---
[void]$ws.Range("E3").Select()
